$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Total Electrical Load (kW)" column (C) values
$ws.Range("C2").Value = 228.15
$ws.Range("C3").Value = 178.15
$ws.Range("C4").Value = 170.79
$ws.Range("C5").Value = 181.12
$ws.Range("C6").Value = 187.57
$ws.Range("C7").Value = 193.32
$ws.Range("C8").Value = 197.26
$ws.Range("C9").Value = 200.94
$ws.Range("C10").Value = 200.19
$ws.Range("C11").Value = 277.14999999999998
$ws.Range("C12").Value = 371.33
$ws.Range("C13").Value = 390.9
$ws.Range("C14").Value = 380.31
$ws.Range("C15").Value = 365.21
$ws.Range("C16").Value = 349.28
$ws.Range("C17").Value = 337.52
$ws.Range("C18").Value = 329.16
$ws.Range("C19").Value = 320.04000000000002
$ws.Range("C20").Value = 310.94
$ws.Range("C21").Value = 304.92
$ws.Range("C22").Value = 307.77999999999997
$ws.Range("C23").Value = 313.06
$ws.Range("C24").Value = 315.16000000000003
$ws.Range("C25").Value = 294.06

# Move selection to H17 (matches saved selection state in the diff)
$ws.Range("H17").Select() | Out-Null
